# feat: add 2022-Q4 data
#
# The workbook had a single quarterly sheet named "2022-Q3" plus a "总计"
# (totals) summary sheet. This adds a new "2022-Q4" quarter:
#   - the existing "2022-Q3" data sheet becomes the new "2022-Q4" sheet
#     (its values are refreshed with the Q4 numbers) while a duplicate of
#     its original content is kept as a brand-new "2022-Q3" sheet, inserted
#     right after it
#   - the "总计" summary sheet's existing row is relabeled "2022-Q4" and a
#     new row repeating the old "2022-Q3" totals is appended below it

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "总计" (totals) sheet: relabel the existing data row as 2022-Q4 and add
# a new row underneath for the (now historical) 2022-Q3 totals.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"

# New row 3 reuses row 2's A-column formatting (bold/bordered style).
$wsTotal.Cells.Item(2, 1).Copy()
$wsTotal.Cells.Item(3, 1).PasteSpecial(-4122)
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(3, 3).Value = 1
$wsTotal.Cells.Item(3, 4).Value = 0.02

# ---------------------------------------------------------------------
# Quarter data sheet: duplicate the current "2022-Q3" sheet (so its old
# data survives unchanged as the new "2022-Q3" tab right after it), then
# rename the original to "2022-Q4" and refresh it with the new numbers.
# ---------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item(2)
$wsQ.Copy([System.Reflection.Missing]::Value, $wsQ)
$wsQ.Name = "2022-Q4"

$wsQOld = $wb.Worksheets.Item(3)
$wsQOld.Name = "2022-Q3"

$wsQ.Cells.Item(2, 4).Value = "'0.67"
$wsQ.Cells.Item(2, 5).Value = "'91.81"
$wsQ.Cells.Item(2, 6).Value = "'2.62"
$wsQ.Cells.Item(2, 7).Value = "'0.0176"
$wsQ.Cells.Item(2, 8).Value = 8

# The leading-apostrophe assignments above force text storage (matching
# the original text-typed numeric cells) but also stamp a quote-prefix
# style; strip that back off so the cells stay visually unstyled like the
# rest of the row.
$wsQ.Cells.Item(1, 1).Copy()
$wsQ.Range("D2:G2").PasteSpecial(-4122)

# Match this refreshed sheet's header-row / leading-column styling to the
# "总计" summary sheet's look (the newly-published quarter tab picks up the
# same header formatting used elsewhere in the workbook).
$wsTotal.Cells.Item(1, 2).Copy()
$wsQ.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Cells.Item(2, 1).Copy()
$wsQ.Cells.Item(2, 1).PasteSpecial(-4122)

# ...and its page margins (points: 0.75in/0.75in/1in/1in/0.5in/0.5in).
$wsQ.PageSetup.LeftMargin = 54
$wsQ.PageSetup.RightMargin = 54
$wsQ.PageSetup.TopMargin = 72
$wsQ.PageSetup.BottomMargin = 72
$wsQ.PageSetup.HeaderMargin = 36
$wsQ.PageSetup.FooterMargin = 36
